$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "color= red;" -> "color: red; font-weight: bold" typo in the
# two existing Volvo-internal-use banner strings (D29 / D30) ---
$v29 = $ws.Range("D29").Value()
$new29 = $v29 -replace [regex]::Escape('color= red;'), 'color: red; font-weight: bold'
$ws.Range("D29").Value = $new29

$v30 = $ws.Range("D30").Value()
$new30 = $v30 -replace [regex]::Escape('color= red;'), 'color: red; font-weight: bold'
$ws.Range("D30").Value = $new30

# --- Append the new row 31: 서비스 테크매니저 포맨 보증 부품 / CP00004 ---
$ws.Range("A31").Value = "서비스 테크매니저 포맨 보증 부품"
$ws.Range("B31").Value = 45453
$ws.Range("C31").Value = "CP00004"

$d31 = @'
<a style="color: red; font-weight: bold">Volvo Internal Use Only</a>
<br> 서비스센터를 통해 문의가 필요한 경우 아래 주소를 참고합니다. 유선 컨택 포인트는 액셀 파일을 참고하십시오.
<br>
<br><h3>에이치모터스</h3>
<br>테크매니저 김도형: dohyung.kim@hvolvo.com
<br><strong>대치: </strong>포맨 신동선: dongsuen.shin@hvolvo.com / 보증, 부품 유현우: hyunwoo.you@hvolvo.com
<br><strong>성수: </strong>포맨 이규용: kyuyong.lee@hvolvo.com / 보증 유병근: byeonggeun.yoo@hvolvo.com / 부품 김태영: taeyoung.kim@hvolvo.com
<br><strong>분당 서현: </strong>포맨 고화정: hwajeong.ko@hvolvo.com / 보증, 부품 김용철: yongcheol.kim@hvolvo.com
<br><strong>율현: </strong>포맨 손경하: kyoungha.son@hvolvo.com / 보증 김현석: hyunseok.kim@hvolvo.com / 부품 김주용: jooyong.kim@hvolvo.com
<br><strong>수원: </strong>포맨 문병원: byeongwon.mun@hvolvo.com / 보증 강병진: byoungjin.kang@hvolvo.com / 부품 김석규: seokgyu.kim@hvolvo.com
<br><strong>인천: </strong>포맨 최동준: dongjun.choi@hvolvo.com / 보증, 부품 조강래: kangrae.cho@hvolvo.com
<br><strong>대전: </strong>포맨 김성식: sungsik.kim@hvolvo.com / 보증, 부품 오명진: myungjin.oh@hvolvo.com
<br><h3>아이비모터스</h3>
<br>테크매니저 김완규: wankyu.kim@ivymotors.co.kr
<br><strong>광주: </strong>포맨 최인호: inho.choi@ivymotors.co.kr / 보증 이근식: geunsik.lee@ivymotors.co.kr / 부품 이병화: byunghwa.lee@ivymotors.co.kr
<br><strong>전주: </strong>포맨 심석순: seoksoon.sim@ivymotors.co.kr / 보증, 부품 박기태: kitae.park@ivymotors.co.kr
<br><strong>순천: </strong>포맨 김민규: mingyu.kim@ivymotors.co.kr / 보증, 부품 문경철지점장: kyungchul.mun@ivymotors.co.kr
<br><strong>제주: </strong>포맨: 최병희: byunghee.choi@ivymotors.co.kr / 보증, 부품 강범석: bumseok.kang@ivymotors.co.kr
<br><h3>태영모터스</h3>
<br>테크매니저 길형근: hyeonggeun.gil@volvoty.co.kr
<br><strong>대구: </strong>포맨 윤주영: juyoung.yun@volvoty.co.kr / 보증, 부품 강범석: bumseok.kang@ivymotors.co.kr
<br><strong>포항: </strong>포맨, 보증 황도근: dogun.hwang@volvoty.co.kr / 부품 이성준: sungjun.lee@volvoty.co.kr
<br><strong>서대구: </strong>포맨 송경규: kyungkyu.song@volvoty.co.kr / 보증, 부품 이재우: jaewoo.lee@volvoty.co.kr
<br><h3>천하자동차</h3>
<br>테크매니저 황원선: wonsun.hwang@chvolvo.co.kr
<br><strong>영등포 양평: </strong>포맨 정진: jin.jeong@chvolvo.co.kr / 보증 이정국: jeongkuk.lee@chvolvo.co.kr / 부품 김성철: sungchul.kim@chvolvo.co.kr
<br><strong>동대문: </strong>포맨 김영환: younghwan.kim@chvolvo.co.kr / 보증, 부품 김종빈: jongbin.kim@chvolvo.co.kr
<br><strong>의정부: </strong>포맨 권순직: soonjik.kwon@chvolvo.co.kr / 보증, 부품 고명진: myoungjin.ko@chvolvo.co.kr
<br><strong>구리: </strong>포맨 한수봉: subong.han@chvolvo.co.kr / 보증, 부품 이종모: jongmo.lee@chvolvo.co.kr
<br><h3>아이언모터스</h3>
<br>테크매니저 박정우: jungwoo.park@ironmotors.co.kr
<br><strong>창원: </strong>포맨 박수호: suho.park@ironmotors.co.kr / 보증 김청석: chungsuk.kim@ironmotors.co.kr / 부품 황준원: junwon.hwang@ironmotors.co.kr
<br><strong>김해: </strong>포맨, 보증 허준우: junewoo.heo@ironmotors.co.kr / 보증 박현진: hyunjin.park@ironmotors.co.kr
<br><strong>광안: </strong>포맨 김봉준: bongjun.kim@ironmotors.co.kr / 보증 유다경: dagyeong.yu@ironmotors.co.kr / 부품 송화섭: hwaseop.song@ironmotors.co.kr
<br><strong>해운대: </strong>포맨 송병근: byeongguen.song@ironmotors.co.kr / 보증 김유빈: youbin.kim@ironmotors.co.kr / 부품 채원호: wonho.chae@ironmotors.co.kr
<br><h3>아주오토리움</h3>
<br>테크매니저 임윤진: yjlim8@aju.co.kr
<br><strong>일산: </strong>포맨 정상준: sjjung2@aju.co.kr / 보증, 부품 모준석: jsmo@aju.co.kr
<br><strong>안양: </strong>포맨 김인태: itkim@aju.co.kr / 보증, 부품 박영오: yopark@aju.co.kr
<br><strong>영등포 문래: </strong>포맨 최지성: jschoi@aju.co.kr / 보증 최형석: hschoi@aju.co.kr / 부품 주인철: icjoo@aju.co.kr
<br><strong>부천: </strong>포맨 이동원: dwlee@aju.co.kr / 보증, 부품 손건희: ghson@aju.co.kr
<br><h3>코오롱 오토모티브</h3>
<br>테크매니저 손영돈: youngdon_son@kolon.com
<br><strong>송파: </strong>포맨 손영돈: youngdon_son@kolon.com / 보증, 부품 장하늘: haneul_jang@kolon.com
<br><strong>서초: </strong>포맨 정승은: seungeun_jeoung@kolon.com / 보증 이윤우: yunwoo_lee@kolon.com / 부품 박용석: yongsuk_park@kolon.com
<br><strong>원주: </strong>포맨 황수근: sugeun_hwang@kolon.com / 보증, 부품 권순영: sunyoung_kwon@kolon.com
<br><strong>천안: </strong>포맨 마정우: jungwoo_ma@kolon.com / 보증, 부품 길원균: wongyun_kil@kolon.com
<br><strong>울산: </strong>포맨 임기범: kibum_lim@kolon.com / 보증, 부품 김민호: minho_kim2@kolon.com
<br><strong>분당 판교: </strong>포맨 조현목: hyunmok_cho1@kolon.com / 보증 임성한: sunghan_im@kolon.com / 부품 정배길: baegil_jung1@kolon.com
<br><strong>강릉: </strong>포맨 석봉근: bonggeun_seok@kolon.com / 보증, 부품 장준영: juneyoung_jang@kolon.com
<br><strong>서산: </strong>포맨 설종환: jonghwan_seol@kolon.com / 보증, 부품 길동현: donghyun_gil@kolon.com

'@
$ws.Range("D31").Value = $d31
$ws.Cells.Item(31, 4).WrapText = $true
$ws.Rows.Item(31).RowHeight = 409.6

$ws.Range("D31").Select() | Out-Null

Write-Output "done"
